$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "channel group / tx test" calc block added to the right of the
# existing CHANNR/FSCTRL table (columns H:Q), mirroring the structure that
# already exists lower on the sheet (rows 14-19) but for the FREQ2/FREQ1/
# FREQ0 registers, plus a small "group id -> freq (GHz)" helper at H1:J2.
# ---------------------------------------------------------------------------

# Row 1 headers
$ws.Range("H1").Value = "group id"
$ws.Range("I1").Value = 2
$ws.Range("O1").Value = "bin"
$ws.Range("P1").Value = "dec"

# Row 2 - group id -> freq helper + FREQ2 register row
$ws.Range("H2").Value = "freq"
$ws.Range("I2").Formula = "=(2400 + (I1)*8)/1000"
$ws.Range("J2").Value = "GHz"

$ws.Range("L2").Value = "#define FREQ2 "
$ws.Range("M2").Value = "0x0D "
$ws.Range("N2").Value = "//Frequency control word, high byte"
$ws.Range("P2").Value = 92
$ws.Range("O2").Formula = "=DEC2BIN(P2)"
$ws.Range("Q2").Formula = "=DEC2HEX(P2)"

# Row 3 - FREQ1 register row
$ws.Range("L3").Value = "#define FREQ1 "
$ws.Range("M3").Value = "0x0E "
$ws.Range("N3").Value = "//Frequency control word, middle byte"
$ws.Range("P3").Value = 236
$ws.Range("O3").Formula = "=DEC2BIN(P3)"

# Row 4 - FREQ0 register row
$ws.Range("L4").Value = "#define FREQ0 "
$ws.Range("M4").Value = "0x0F "
$ws.Range("N4").Value = "//Frequency control word, low byte"
$ws.Range("P4").Value = 79
$ws.Range("O4").Formula = "=DEC2BIN(P4)"

# Q3:Q4 share one formula (DEC2HEX), entered as a single fill so it is
# stored as a shared formula group, matching Q2 being its own formula.
$ws.Range("Q3:Q4").Formula = "=DEC2HEX(P3)"

# Row 6 - bit-width helper for the new P2:P4 block
$ws.Range("O6").Formula = "=LEN(DEC2BIN(LARGE(P2:P4,1)))"

# Row 7 - bin/dec sub headers
$ws.Range("M7").Value = "bin"
$ws.Range("N7").Value = "dec"

# Row 8 - fxosc
$ws.Range("L8").Value = "fxosc"
$ws.Range("N8").Value = 26
$ws.Range("O8").Value = "MHz"

# Row 9 - base freq
$ws.Range("L9").Value = "base freq"
$ws.Range("M9").Formula = "=CONCAT(O4,O3,O2)"
$ws.Range("N9").Formula = "=P4+P3*2^8+P2*2^16"

# Row 11 - fcarrier (Hz)
$ws.Range("L11").Value = "fcarrier"
$ws.Range("M11").Formula = "=(N8*10^6)/2^16*(N9)"
$ws.Range("N11").Value = "Hz"

# Row 12 - fcarrier (GHz)
$ws.Range("L12").Value = "fcarrier"
$ws.Range("M12").Formula = "=M11/10^9"
$ws.Range("N12").Value = "GHz"

# ---------------------------------------------------------------------------
# Column M width - widened (autofit) after the new numbers were entered.
# ---------------------------------------------------------------------------
[void]$ws.Columns.Item(13).AutoFit()

# ---------------------------------------------------------------------------
# Conditional formatting: keep the existing "D13 > 8" highlight rule, and
# add the same highlight rule on the new O6 bit-width cell. The new rule is
# created first (on the existing rule's dxf) and then the old rule is
# re-added so it ends up as the second / lower priority rule, matching the
# target layout (O6 -> priority 1, D13 -> priority 2).
# ---------------------------------------------------------------------------
$existingRule = $ws.Range("D13").FormatConditions.Item(1)
$existingRule.ModifyAppliesToRange($ws.Range("O6"))

$newRule = $ws.Range("D13").FormatConditions.Add(1, 5, "8")
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# Page setup - explicit portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selection, to match the author's last cursor position.
# ---------------------------------------------------------------------------
[void]$ws.Range("P5").Select()
